# Integrating way to update payment history
# Append a new user row (row 4) to the "User" sheet, including a mailto
# hyperlink on the email cell, mirroring the existing rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$email = "kennedy.ads@gmail.com"

$ws.Range("B4").Value = $email
$ws.Range("C4").Value = 710943451
$ws.Range("D4").Value = "`$2b`$10`$QErIzwxCF2DMPC8bd2TPou7uEnEm.bsZ6lahX1mAufcA5SF.hnaP2"
$ws.Range("E4").Value = $true

# Add the mailto hyperlink on the new email cell (matches B2/B3 pattern),
# then strip the auto-applied "Hyperlink" style so formatting matches the
# existing (unstyled) rows.
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:" + $email)
$ws.Range("B4").Style = "Normal"
